# Add a new "ATS Competitors" worksheet at the end of the workbook (after
# "ATS History") containing the Grok-verified ATS resume tool competitor
# analysis plus the MatchForge differentiation / honest-assessment tables.

$wb = $excel.ActiveWorkbook

# Insert the new sheet directly after the current last sheet ("ATS History")
# so it lands as the final, 8th tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ATS Competitors"

# --- Section 1: Competitor pricing table -------------------------------
$ws.Range("A1").Value = "ATS RESUME TOOL COMPETITORS (Grok-verified Jan 2026)"

$ws.Range("A3").Value = "Tool"
$ws.Range("B3").Value = "Free Tier"
$ws.Range("C3").Value = "Paid Monthly"
$ws.Range("D3").Value = "Key Feature"

$ws.Range("A4").Value = "Jobscan"
$ws.Range("B4").Value = "Limited scans"
$ws.Range("C4").Value = "$49.95/mo"
$ws.Range("D4").Value = "Market leader, job matching"

$ws.Range("A5").Value = "SkillSyncer"
$ws.Range("B5").Value = "Limited"
$ws.Range("C5").Value = "$19.99-39.99/mo"
$ws.Range("D5").Value = "AI auto-optimization"

$ws.Range("A6").Value = "Resume Worded"
$ws.Range("B6").Value = "Basic checks"
$ws.Range("C6").Value = "$8/mo"
$ws.Range("D6").Value = "25+ criteria checks"

$ws.Range("A7").Value = "Enhancv"
$ws.Range("B7").Value = "Basic builder"
$ws.Range("C7").Value = "~$14/mo"
$ws.Range("D7").Value = "Resume builder + ATS"

$ws.Range("A8").Value = "ResyMatch.io"
$ws.Range("B8").Value = "Limited"
$ws.Range("C8").Value = "$8/mo"
$ws.Range("D8").Value = "Free tier from Cultivated Culture"

$ws.Range("A9").Value = "MatchForge"
$ws.Range("B9").Value = "Full demo"
$ws.Range("C9").Value = "$9-29/mo (target)"
$ws.Range("D9").Value = "ATS + Matching + Coaching"

# --- Section 2: True differentiation matrix -----------------------------
$ws.Range("A11").Value = "MATCHFORGE TRUE DIFFERENTIATION"

$ws.Range("A12").Value = "Feature"
$ws.Range("B12").Value = "Jobscan"
$ws.Range("C12").Value = "SkillSyncer"
$ws.Range("D12").Value = "MatchForge"

$ws.Range("A13").Value = "ATS checking"
$ws.Range("B13").Value = "Yes"
$ws.Range("C13").Value = "Yes"
$ws.Range("D13").Value = "Yes"

$ws.Range("A14").Value = "Transparent 6-factor matching"
$ws.Range("B14").Value = "No"
$ws.Range("C14").Value = "No"
$ws.Range("D14").Value = "YES"

$ws.Range("A15").Value = "Human coaching"
$ws.Range("B15").Value = "No"
$ws.Range("C15").Value = "No"
$ws.Range("D15").Value = "YES"

$ws.Range("A16").Value = "Outcome tracking"
$ws.Range("B16").Value = "No"
$ws.Range("C16").Value = "No"
$ws.Range("D16").Value = "YES"

$ws.Range("A17").Value = "Feedback loop validation"
$ws.Range("B17").Value = "No"
$ws.Range("C17").Value = "No"
$ws.Range("D17").Value = "YES"

$ws.Range("A18").Value = "Full stack (ATS+Match+Coach)"
$ws.Range("B18").Value = "No"
$ws.Range("C18").Value = "No"
$ws.Range("D18").Value = "YES"

# --- Section 3: Honest assessment ---------------------------------------
$ws.Range("A20").Value = "HONEST ASSESSMENT"
$ws.Range("A21").Value = "ATS checking alone is NOT unique. The COMBINATION is the differentiator."

# Land the selection back on A1, matching a freshly-opened sheet.
$ws.Range("A1").Select() | Out-Null
